$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rotate F:V for rows 69-72 (old row72 data -> row69; old rows 69-71 shift down to 70-72) ---
$ws.Range("F69").Value = 'Zielona Gora'
$ws.Range("G69").Value = 2
$ws.Range("H69").Value = 'Bytom Odrzanski'
$ws.Range("I69").Value = 2
$ws.Range("J69").Value = 2.3
$ws.Range("K69").Value = '29/09/2023 03:12'
$ws.Range("L69").Value = 2.05
$ws.Range("M69").Value = '30/09/2023 15:51'
$ws.Range("N69").Value = 3.21
$ws.Range("O69").Value = '29/09/2023 03:12'
$ws.Range("P69").Value = 3.55
$ws.Range("Q69").Value = '30/09/2023 15:51'
$ws.Range("R69").Value = 2.5
$ws.Range("S69").Value = '29/09/2023 03:12'
$ws.Range("T69").Value = 2.97
$ws.Range("U69").Value = '30/09/2023 15:51'
$ws.Range("V69").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-iii/zielona-gora-bytom-odrzanski/nwCsPr37/'

$ws.Range("F70").Value = 'Carina Gubin'
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 'Sleza Wroclaw'
$ws.Range("I70").Value = 4
$ws.Range("J70").Value = 2.43
$ws.Range("K70").Value = '29/09/2023 03:12'
$ws.Range("L70").Value = 2.47
$ws.Range("M70").Value = '30/09/2023 15:58'
$ws.Range("N70").Value = 3.32
$ws.Range("O70").Value = '29/09/2023 03:12'
$ws.Range("P70").Value = 3.68
$ws.Range("Q70").Value = '30/09/2023 15:58'
$ws.Range("R70").Value = 2.31
$ws.Range("S70").Value = '29/09/2023 03:12'
$ws.Range("T70").Value = 2.34
$ws.Range("U70").Value = '30/09/2023 15:58'
$ws.Range("V70").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-iii/carina-gubin-sleza-wroclaw/z13aEn6n/'

$ws.Range("F71").Value = 'Pawlowice'
$ws.Range("G71").Value = 3
$ws.Range("H71").Value = 'Gornik Zabrze II'
$ws.Range("I71").Value = 2
$ws.Range("J71").Value = 1.91
$ws.Range("K71").Value = '29/09/2023 03:12'
$ws.Range("L71").Value = 1.72
$ws.Range("M71").Value = '30/09/2023 15:07'
$ws.Range("N71").Value = 3.45
$ws.Range("O71").Value = '29/09/2023 03:12'
$ws.Range("P71").Value = 3.9
$ws.Range("Q71").Value = '30/09/2023 15:07'
$ws.Range("R71").Value = 3
$ws.Range("S71").Value = '29/09/2023 03:12'
$ws.Range("T71").Value = 3.74
$ws.Range("U71").Value = '30/09/2023 15:07'
$ws.Range("V71").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-iii/pniowek-pawlowice-gornik-zabrze/vD8ZQ4Yf/'

$ws.Range("F72").Value = 'Stilon Gorzow'
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = 'Goczalkowice Zdroj'
$ws.Range("I72").Value = 1
$ws.Range("J72").Value = 2.53
$ws.Range("K72").Value = '29/09/2023 03:12'
$ws.Range("L72").Value = 2.9
$ws.Range("M72").Value = '30/09/2023 15:58'
$ws.Range("N72").Value = 3.16
$ws.Range("O72").Value = '29/09/2023 03:12'
$ws.Range("P72").Value = 3.45
$ws.Range("Q72").Value = '30/09/2023 15:58'
$ws.Range("R72").Value = 2.3
$ws.Range("S72").Value = '29/09/2023 03:12'
$ws.Range("T72").Value = 2.12
$ws.Range("U72").Value = '30/09/2023 15:58'
$ws.Range("V72").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-iii/stilon-gorzow-goczalkowice-zdroj/QTAoO2ID/'

# --- Append new row 97 ---
$ws.Cells.Item(96,1).Copy()
$ws.Cells.Item(97,1).PasteSpecial(-4122)
$ws.Cells.Item(96,5).Copy()
$ws.Cells.Item(97,5).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A97").Value = 96
$ws.Range("B97").Value = 'poland'
$ws.Range("C97").Value = 'iii-liga-group-iii'
$ws.Range("D97").Value = '2023-2024'
$ws.Range("E97").Value = 45224.625
$ws.Range("F97").Value = 'Zielona Gora'
$ws.Range("G97").Value = 2
$ws.Range("H97").Value = 'Gornik Zabrze II'
$ws.Range("I97").Value = 1
$ws.Range("J97").Value = 2.05
$ws.Range("K97").Value = '24/10/2023 02:12'
$ws.Range("L97").Value = 2.08
$ws.Range("M97").Value = '25/10/2023 14:59'
$ws.Range("N97").Value = 3.43
$ws.Range("O97").Value = '24/10/2023 02:12'
$ws.Range("P97").Value = 3.2
$ws.Range("Q97").Value = '25/10/2023 14:59'
$ws.Range("R97").Value = 2.71
$ws.Range("S97").Value = '24/10/2023 02:12'
$ws.Range("T97").Value = 2.9
$ws.Range("U97").Value = '25/10/2023 14:59'
$ws.Range("V97").Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-iii/zielona-gora-gornik-zabrze/tY0oONZA/'